$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.253.84'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '3.549.04'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.05'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').Value = '3.547.37'
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.489'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.08%  '
$ws.Range('E10').Value = '  -0.77%  '
$ws.Range('E11').Value = '  -2.93%  '
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').Value = '4.152.98'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('E14').Value = '  -0.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '30.03'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('D16').Value = '3.556.47'
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('D17').Value = '66.353.33'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.46'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.93%  '
$ws.Range('E20').Value = '  -0.81%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.79'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '430.44'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.609'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.55'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Value = '3.693.45'
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('E27').Value = '  -0.94%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.50'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.31%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.93'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.44%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.10'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '25.38'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('B33').Value = 'RenzoRestakedETH'
$ws.Range('C33').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D33').Value = '3.544.05'
$ws.Range('E33').Value = '  +0.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.44'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.71%  '
$ws.Range('E35').Value = '  -5.75%  '
$ws.Range('B36').Value = 'USDe'
$ws.Range('C36').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.81'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.88%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.72'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.64%  '
$ws.Range('E39').Value = '  -0.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '175.97'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.33%  '
$ws.Range('E41').Value = '  -1.50%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.18'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.886'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.78%  '
$ws.Range('E44').Value = '  +1.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '45.94'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.50'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.56%  '
$ws.Range('E48').Value = '  -1.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.11'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.12'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.19'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.37%  '
